# edit.ps1 - apply the crypto-list data refresh described in the commit
# "Updated cryptos list on Wed Mar  8 21:53:33 UTC 2023 with GitHub Actions".
#
# Rows 2-51 on Sheet1 hold one coin per row: B=Coin name, C=Link, D=Price,
# E=Volume(1h). This refresh updates nearly every Price/Volume cell with the
# newly scraped figures; two pairs of rows (15/16 and 45/46) also swap which
# coin occupies which row (ranking reshuffled), so their Coin/Link/Price/Volume
# are all rewritten together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "21.967.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.48%  "

# Row 3: Ethereum -> Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.552.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.03%  "

# Row 4: TetherUSD -> TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.28%  "

# Row 5: USDC -> USDC
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.08%  "

# Row 6: BNB -> BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.17%  "

# Row 7: XRP -> XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3922"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.12%  "

# Row 8: Cardano -> Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3225"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.65%  "

# Row 9: OKB -> OKB
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.40"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.81%  "

# Row 10: Dogecoin -> Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07185"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.39%  "

# Row 11: Polygon -> Polygon
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.070"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.89%  "

# Row 12: BinanceUSD -> BinanceUSD
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.29%  "

# Row 13: Polkadot -> Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.655"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.35%  "

# Row 14: Solana -> Solana
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.19%  "

# Row 15: Chainlink -> WrappedEther
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.563.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.04%  "

# Row 16: WrappedEther -> Chainlink
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.630"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.80%  "

# Row 17: ShibaInu -> ShibaInu
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001109"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.09%  "

# Row 18: TRON -> TRON
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06595"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.49%  "

# Row 19: Litecoin -> Litecoin
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.81%  "

# Row 20: Dai -> Dai
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.03%  "

# Row 21: Uniswap -> Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.218"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.47%  "

# Row 22: Avalanche -> Avalanche
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.82%  "

# Row 23: Cosmos -> Cosmos
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.75%  "

# Row 24: WrappedBTC -> WrappedBTC
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "21.987.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.44%  "

# Row 25: Toncoin -> Toncoin
$ws.Range("E25").Value = "  +2.65%  "

# Row 26: LidoDAOToken -> LidoDAOToken
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.387"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.19%  "

# Row 27: Monero -> Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "147.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.88%  "

# Row 28: EthereumClassic -> EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.92%  "

# Row 29: HuobiToken -> HuobiToken
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.860"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.31%  "

# Row 30: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.736.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.93%  "

# Row 31: BitcoinCash -> BitcoinCash
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.01%  "

# Row 32: ImmutableX -> ImmutableX
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9779"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.31%  "

# Row 33: Filecoin -> Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.880"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.04%  "

# Row 34: Stellar -> Stellar
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08290"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.79%  "

# Row 35: FraxShare -> FraxShare
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.107"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.40%  "

# Row 36: WEMIXTOKEN -> WEMIXTOKEN
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.610"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -15.24%  "

# Row 37: VeChain -> VeChain
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02253"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.71%  "

# Row 38: InternetComputer(DFINITY) -> InternetComputer(DFINITY)
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.094"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.24%  "

# Row 39: Hedera -> Hedera
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05994"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.85%  "

# Row 40: TrustWalletToken -> TrustWalletToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.204"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.31%  "

# Row 41: Algorand -> Algorand
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2051"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.57%  "

# Row 42: Frax -> Frax
$ws.Range("E42").Value = "  +0.00%  "

# Row 43: Aptos -> Aptos
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.92%  "

# Row 44: TheSandbox -> TheSandbox
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5777"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.72%  "

# Row 45: EnergySwap -> PancakeSwap
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.749"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.56%  "

# Row 46: PancakeSwap -> EnergySwap
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.55%  "

# Row 47: Decentraland -> Decentraland
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5545"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.58%  "

# Row 48: Quant -> Quant
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "117.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.94%  "

# Row 49: NEARProtocol -> NEARProtocol
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.873"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.80%  "

# Row 50: EOS -> EOS
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.133"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.40%  "

# Row 51: Cronos -> Cronos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06810"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.97%  "
